$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2:C51").ClearContents()
$ws.Range("B2").Value = "_tejgct_r09gstcp"
$ws.Range("C2").Value = 0.0222877008616479
$ws.Range("B3").Value = "_tejgkft_redr"
$ws.Range("C3").Value = 0.01876119338234825
$ws.Range("B4").Value = "_tejgge_r09ct06acanf"
$ws.Range("C4").Value = 0.01270238711430756
$ws.Range("B5").Value = "_devppimtotfun_f1trans"
$ws.Range("C5").Value = 0.008379448440144543
$ws.Range("B6").Value = "_tejgfun_f5ct05prots"
$ws.Range("C6").Value = 0.007266205920448248
$ws.Range("B7").Value = "_tejgfun_f5ct06amb"
$ws.Range("C7").Value = 0.006509051985726964
$ws.Range("B8").Value = "devppimfun_f5r07ct05salud"
$ws.Range("C8").Value = 0.0057942148843925
$ws.Range("B9").Value = "tejgfun_f5ct06viv"
$ws.Range("C9").Value = 0.005213637077382253
$ws.Range("B10").Value = "tejgtotfun_f2prots"
$ws.Range("C10").Value = 0.004362137750170723
$ws.Range("B11").Value = "_tejgfun_f5ct06opseg"
$ws.Range("C11").Value = 0.003752831168538787
$ws.Range("B12").Value = "_tejgfun_f5r18ct05pgrco"
$ws.Range("C12").Value = 0.003543889931279804
$ws.Range("B13").Value = "_tejgtotfun_f2opsegpc"
$ws.Range("C13").Value = 0.002667236308274454
$ws.Range("B14").Value = "pimgfun_f5r18ct05trans"
$ws.Range("C14").Value = 0.002618618194689319
$ws.Range("B15").Value = "_tejgfun_f5r08ct05prots"
$ws.Range("C15").Value = 0.002402956011004317
$ws.Range("B16").Value = "dfgdevpiagfun_f5ct05sanpc"
$ws.Range("C16").Value = 0.002328877006073631
$ws.Range("B17").Value = "devppimtotfun_f5r07salud"
$ws.Range("C17").Value = 0.002327153653263374
$ws.Range("B18").Value = "dfgpimpiafun_f1ct05prots"
$ws.Range("C18").Value = 0.002320550390594772
$ws.Range("B19").Value = "_tejgtotfun_f5r08pgrco"
$ws.Range("C19").Value = 0.002167362117065952
$ws.Range("B20").Value = "tejgfun_f5ct05trans"
$ws.Range("C20").Value = 0.002118668341594504
$ws.Range("B21").Value = "_pimgfun_f5ct06opsegpc"
$ws.Range("C21").Value = 0.0020550167627146
$ws.Range("B22").Value = "devppimfun_f5r07ct05agro"
$ws.Range("C22").Value = 0.002020945930915277
$ws.Range("B23").Value = "_devppimfun_f1ct05trans"
$ws.Range("C23").Value = 0.002011812909351638
$ws.Range("B24").Value = "tejgfun_f5r08ct05ambpc"
$ws.Range("C24").Value = 0.001939748180054954
$ws.Range("B25").Value = "_tejgtotfun_f5r18prots"
$ws.Range("C25").Value = 0.001865224340092256
$ws.Range("B26").Value = "tejgfun_f2ct05ambpc"
$ws.Range("C26").Value = 0.001783193650746896
$ws.Range("B27").Value = "tejgtotfun_f5r18opseg"
$ws.Range("C27").Value = 0.001696893484526554
$ws.Range("B28").Value = "_pimgfun_f5r18ct06opseg"
$ws.Range("C28").Value = 0.001647335556545604
$ws.Range("B29").Value = "dfgpimpiatotfun_f1prots"
$ws.Range("C29").Value = 0.001623508976675453
$ws.Range("B30").Value = "_tejgtotfun_f2opseg"
$ws.Range("C30").Value = 0.001610369828400165
$ws.Range("B31").Value = "_tejgtotfun_f5amb"
$ws.Range("C31").Value = 0.001606209415203525
$ws.Range("B32").Value = "devppimfun_f5ct06viv"
$ws.Range("C32").Value = 0.001603637616308873
$ws.Range("B33").Value = "_tejgfun_f2ct05pgrco"
$ws.Range("C33").Value = 0.00160282475040108
$ws.Range("B34").Value = "tejgfun_f2ct05prots"
$ws.Range("C34").Value = 0.001575118861645958
$ws.Range("B35").Value = "_tejgfun_f5r18ct05prots"
$ws.Range("C35").Value = 0.001522255742288848
$ws.Range("B36").Value = "pimgfun_f1ct06san"
$ws.Range("C36").Value = 0.001516097295823236
$ws.Range("B37").Value = "_tejgfun_f5ct06opsegpc"
$ws.Range("C37").Value = 0.001507765009685699
$ws.Range("B38").Value = "dfgdevpiagfun_f5r18ct05transpc"
$ws.Range("C38").Value = 0.001450216865727935
$ws.Range("B39").Value = "tejgct_r07gstcrpc"
$ws.Range("C39").Value = 0.001403685900253843
$ws.Range("B40").Value = "devppimfun_f5r07ct05cydep"
$ws.Range("C40").Value = 0.001353879120617375
$ws.Range("B41").Value = "devppimfun_f5r18ct05agro"
$ws.Range("C41").Value = 0.001344243803332496
$ws.Range("B42").Value = "_devppimfun_f5r07ct05trans"
$ws.Range("C42").Value = 0.001343546795418223
$ws.Range("B43").Value = "devppimtotfun_f5r07edu"
$ws.Range("C43").Value = 0.001342745141028111
$ws.Range("B44").Value = "tejgfun_f5ct05prots"
$ws.Range("C44").Value = 0.001332322578335413
$ws.Range("B45").Value = "devppimtotfun_f5turi"
$ws.Range("C45").Value = 0.001316269411181294
$ws.Range("B46").Value = "_tejgtotfun_f2protspc"
$ws.Range("C46").Value = 0.001299328638028498
$ws.Range("B47").Value = "_tejgtotfun_f5r18amb"
$ws.Range("C47").Value = 0.001289593658318534
$ws.Range("B48").Value = "_devppimfun_f5r07ct05cydep"
$ws.Range("C48").Value = 0.001224440924213341
$ws.Range("B49").Value = "_tejgge_r00ct05biser"
$ws.Range("C49").Value = 0.001219903169780348
$ws.Range("B50").Value = "_devppimtotfun_f5viv"
$ws.Range("C50").Value = 0.00120469137047468
$ws.Range("B51").Value = "_tejgge_r08ct05pobso"
$ws.Range("C51").Value = 0.001202337329518565
